$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F "想去人数" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 73
$ws1.Range("F7").Value = 87
$ws1.Range("F13").Value = 2378
$ws1.Range("F16").Value = 524
$ws1.Range("F17").Value = 547
$ws1.Range("F18").Value = 163
$ws1.Range("F19").Value = 83
$ws1.Range("F20").Value = 48
$ws1.Range("F22").Value = 1877
$ws1.Range("F23").Value = 4012
$ws1.Range("F28").Value = 2092
$ws1.Range("F34").Value = 420
$ws1.Range("F36").Value = 694
$ws1.Range("F37").Value = 437
$ws1.Range("F38").Value = 415

# Sheet "演出" (sheet2) - column F "想去人数" update
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 34

# Sheet "全部类型" (sheet4) - column F "想去人数" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 73
$ws4.Range("F7").Value = 87
$ws4.Range("F13").Value = 2378
$ws4.Range("F15").Value = 34
$ws4.Range("F17").Value = 524
$ws4.Range("F18").Value = 547
$ws4.Range("F19").Value = 163
$ws4.Range("F20").Value = 83
$ws4.Range("F21").Value = 48
$ws4.Range("F23").Value = 1877
$ws4.Range("F24").Value = 4012
$ws4.Range("F29").Value = 2092
$ws4.Range("F35").Value = 420
$ws4.Range("F37").Value = 694
$ws4.Range("F38").Value = 437
$ws4.Range("F39").Value = 415
